$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "Caption" header for the Images table (row 5), new column G
$ws.Range("G5").Value = "Caption"

# Set width for the new column G to match the diff (stored width="24")
$ws.Columns.Item(7).ColumnWidth = 23.1

# Update selection to match the post-edit state (G6)
$ws.Range("G6").Select()
